$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 ---
$A3 = @'
첫구매는 반값다~딜 
'@
$ws.Range("A3").Value = $A3
$B3 = @'
첫구매는 세 개 골라 다~반값 + 무료배송
'@
$ws.Range("B3").Value = $B3
$C3 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000004059&domainSiteNo=7018
'@
$ws.Range("C3").Value = $C3
$D3 = @'
반값다딜 - CJ편 (2/1~8)
'@
$ws.Range("D3").Value = $D3
$ws.Range("E3").NumberFormat = "@"
$E3 = @'
2023-02-01
'@
$ws.Range("E3").Value = $E3
$ws.Range("F3").NumberFormat = "@"
$F3 = @'
2023-02-08
'@
$ws.Range("F3").Value = $F3
$G3 = @'
['이벤트/쿠폰 > 반값다딜 - CJ편 (2/1~8)', '스마일클럽', '첫구매 전용 프로모션', '50% 할인 쿠폰, 무료배송', '첫구매는 세 개 골라 반값!', '* 첫구매 고객은 생애 첫 구매 고객과 최근 1년 이내 구매 이력이 없는 고객을 포함합니다.', 'cj상품 50% 할인 쿠폰 & 첫구매 무료배송 쿠폰', ' 최대 할인 금액 각 5천원 ', '         쿠폰 발급 기간', '         쿠폰 발급 대상', '         2022년 1월 31일 이후 이마트몰, 트레이더스 쓱배송/점포택배 및 새벽배송 구매 이력이 없는 고객 한정', '         상품 할인 쿠폰 적용 방법', '         본 페이지 내 전시된 CJ브랜드 상품 중 3개를 골라 결제 단계에서 쿠폰 적용 (상품당 1개 수량 적용)', '         20,000원 이상 구매시 사용 가능', '         무료배송 쿠폰은 2만원 이상 구매시 사용가능하며, 반값쿠폰과 함께 발급됩니다.', '첫구매 쿠폰 모두 한번에 받기', '쿠폰 사용 전 꼭 확인하세요!', '       쿠폰 사용 기간', '       상품할인쿠폰/무료배송 쿠폰 : 2023년 2/1(수) ~ 2/8(수) 행사 기간 내 사용', '       쿠폰 사용 조건', '       상품할인 쿠폰 : 총 3장, 할인 적용 가능 상품들에 한해 50% 할인(상품당 최대할인금액 적용)', '       무료배송 쿠폰 : 총 1장, 이마트몰 쓱배송/새벽배송 상품 2만원 이상 구매시 무료배송', '       쿠폰 발급 대상', '       쿠폰 적용 대상 상품', '       상품당 1개 수량에 쿠폰 적용 가능하며, 동일한 상품 2개 구매시에도 1개에만 적용 가능합니다.']
'@
$ws.Range("G3").Value = $G3

# --- Row 4 ---
$A4 = @'
더 강력해진 SSG.COM 삼성카드
'@
$ws.Range("A4").Value = $A4
$B4 = @'
SSG MONEY 최대 15% 적립 + 스마일클럽 월이용료 할인
'@
$ws.Range("B4").Value = $B4
$C4 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000002385&recruitmentPath=SSG
'@
$ws.Range("C4").Value = $C4
$D4 = @'
SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지
'@
$ws.Range("D4").Value = $D4
$ws.Range("E4").NumberFormat = "@"
$E4 = @'
2022-10-26
'@
$ws.Range("E4").Value = $E4
$ws.Range("F4").NumberFormat = "@"
$F4 = @'
2025-10-25
'@
$ws.Range("F4").Value = $F4
$G4 = @'
['이벤트/쿠폰 > SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지', '스마일클럽', 'SSG.COM삼성카드 리뉴얼 이벤트 안내 페이지', 'SSG머니 최대 15% 적립', '스마일클럽 월 이용료 3,900원 결제일 할인', '최대 15% SSG머니 적립 중 5%는 SSG.COM에서 제공하는 서비스로 자세한 내용은 SSG.COM 확인바람.', '01. SSG MONEY 최대 15% 적립!', '       이용실적 관계없이 적립 : 1 % + 이용실적 충족 시 적립 4% = 적립 가능한 최대 SSG머니 5%', '       카드혜택 + 스마일클럽 혜택', '       이용실적 관계없이 적립 1% + 이용실적 충족 시 적립 9% + 쓱배송/쓱배송 traders/새벽배송 상품 구매 시 5% = 적립 가능한 최대 SSG머니 15%', '최대 15% SSG MONEY 적립', 'SSG.COM 삼성카드 최대 10% 적립+스마일클럽 5% 적립', 'SSG.COM 삼성카드 최대 10% 적립(1% 적립+9% 추가 적립)', '1% 적립(전월 실적 조건 및 적립한도 없음)', '9% 추가 적립', '삼성카드 할인이 적용된 일시불 및 할부 이용금액은 제외됩니다.', '스마일클럽 5% 적립', '쓱배송/쓱배송 traders/새벽배송 상품 결제건에 한해 혜택을 받을 수 있습니다.', '02. 스마일클럽 월이용료 매월 3,900원 결제일 할인', '15,000원만 담아도 쓱 새벽배송 무료배송(SSG가입시)', '1쓱 새벽배송 트레이더 - SSG머니 최대 5% 적립', '장바구니 최대 10% 할인쿠폰', '스마일클럽 단독 혜택 - 스타벅스 상품 전용 딜', '매월 4장씩 최대 12% 할인쿠폰', '스마일배송 1만 5천원 이상 무제한 무료배송', '스마일배송 상품 스마일캐시 1% 적립', '카드 할인 혜택 자세히보기 (레이어팝업 열기)', '2022.10.26부터, SSG.COM 삼성카드 발급 시 스마일클럽에 자동 가입됩니다.', 'G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우 혜택을 받을 수 없습니다.', 'SSG.COM 삼성카드로 스마일클럽 월 이용료(3,900원) 정기결제 시 혜택이 제공됩니다.(월 1회)', 'SSG.COM을 통한 스마일클럽 가입건에 한해 혜택이 제공됩니다.', '결제금액이 할인금액보다 적을 경우, 결제금액만큼 할인이 적용됩니다.(결제금액이 없는 경우 할인 대상에서 제외)', '01. 5만 5천원 이상 결제 시 사용 가능한 5만원 할인쿠폰 제공', '50,000원 할인쿠폰 - SSGPAY 바로결제 이용 시 사용 가능', '쿠폰 발급기간 : 2023.02.01 ~ 2023.02.28', '쿠폰 사용기간 : 2023.02.01 ~ 2023.02.28', 'SSGPAY 바로결제 이용 이력이 없는 회원', 'SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제해야 쿠폰을 사용할 수 있습니다.', '배송비 등을 제외한 최종 결제금액이 55,000원 이상이여야 쿠폰을 사용할 수 있습니다.', 'SSGPAY 바로결제에 등록된 SSG.COM 삼성카드로 결제하셔야 쿠폰을 사용할 수 있습니다.', '쿠폰은 통합 회원 본인 명의의 SSG.COM 삼성카드로 결제 시 사용 가능하며, 다른 부정적인 방법으로 사용한 경우에는 주문이 취소될 수 있습니다.', '쿠폰은 결제 화면에서 자동으로 적용됩니다.', '결제 화면에서 쿠폰 변경을 원할 경우 ‘쿠폰선택’을 눌러주세요.', '02. 삼성카드 결제일할인 행사일에는 최대 5% 즉시할인', '[카드혜택가 > SSG.COM 삼성카드 행사 더보기] 내 즉시할인가가 표시된 상품 결제시 적용(상품별 할인한도 30,000원)', '삼성카드 결제일할인 행사일에만 적용됩니다.', '삼성카드 할인이 적용된 결제건은 SSG머니 적립 대상에서 제외됩니다.', '03. 스마일클럽 월 이용료 결제 카드로 SSG.COM 삼성카드 등록 시 SSG머니 3,900원 적립', '행사기간 동안 SSG.COM 삼성카드를 통해 스마일클럽 자동 가입 시 SSG머니 3,900원 즉시 적립', '적립시점까지 스마일클럽 월 이용료 정기결제 수단에 SSG.COM 삼성카드를 등록해야 혜택을 받을 수 있습니다.', '본 상품 발급 및 SSGPAY앱 설치시 바로결제에 자동 등록 됩니다.', 'SSG MONEY는 매월 1일~말일까지 매출전표가 접수된 금액에 대해 다음달 25일 SSG.COM 계정으로 적립됩니다.', '적립된 SSG MONEY는 SSGPAY 회원가입 후 조회 및 사용 가능합니다. 단, SSG.COM 에서는 SSG.COM만 가입해도 사용 가능합니다.', '연체이자율 : 회원별/이용상품별 정상이자율+3.0%p(최고 연 20.0%)', '이미 SSG.COM 삼성카드를 가지고 계시네요!스마일클럽 가입하고 모든 혜택 누리세요', '이미 SSG.COM 삼성카드를 가지고 계시네요!G마켓 또는 옥션을 통해 스마일클럽에 가입한 경우, 월 이용료 결제일할인 혜택을 받을 수 없습니다.']
'@
$ws.Range("G4").Value = $G4

# --- Row 5 ---
$A5 = @'
최대 5만원 혜택
'@
$ws.Range("A5").Value = $A5
$G5 = @'
['이벤트/쿠폰 > SSG.COM카드 Edition 2 이벤트 안내 페이지', '스마일클럽', 'SSG.COM카드 Edition 2 이벤트 안내 페이지', 'SSG.COM카드 Edition2는 SSG.COM에서 최대 5만원 혜택 + 쓸 때마다 최대 12% 적립', '이달의 혜택 01. SSG.COM카드 Edition2 첫결제 시 1만원 쿠폰 할인 바로보기', '혜택 01. SSG머니 최대 12% 적립 바로보기', '혜택 02. 매월 스마일클럽 가입비 3,900원 지원 바로보기', '이달의 혜택 01. SSG.COM에서 SSG.COM카드 Edition2 첫 결제 시 1만원 할인', '      직전 6개월간 (2022년 8월 1일 부터 2023년 1월 31일) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한함', '첫 결제 쿠폰 이용 방법', '② SSGPAY에 등록된 SSG.COM카드 Edition2로', '③ SSG.COM에서 기간 내 1만 1천원 이상 첫 결제 시 1만원 쿠폰 할인', '      본 혜택은 SSG.COM카드 및 SSG.COM카드 Edition2 로 SSG.COM(이마트몰, 신세계몰, 신세계백화점몰 등)에서 직전 6개월간 (2022년 7월 1일 부터 년 12월 31일 까지) SSG.COM카드 및 SSG.COM카드 Edition2로 결제 이력이 없고 & SSGPAY에 등록된 현대카드가 없는 회원에 한해 제공됨', "      본 혜택은 결제 시 [결제방법 > SSGPAY카드] 내 SSG.COM카드 Edition2 선택 시, '카드할인 최적가' 추천에 의해 할인 금액이 자동 적용됨. 단, 1만 1천원 이상 첫 결제 시 적용)", '다운받은 쿠폰은 SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.', '쿠폰 사용하러 가기', 'SSGPAY에 등록된 SSG.COM카드 Edition2로 SSG.COM에서 바로 결제 시 사용 가능합니다.', '4. SSG.COM에서 기간 내 2만원 이상 첫 결제 시 제공 *쿠폰할인, 적립금 사용(SSG머니, 신세계포인트 등), 상품권 등 선할인 및 일부상품(상품권 등 현금성 상품, 무형서비스 상품, 초특가 상품, 특정 브랜드 상품 등)을 제외한 카드 최종 결제금액이 2만원 이상인 경우에 한해 혜택 적용', '혜택 01. 장 볼 때마다 SSG머니 최대 12% 적립', 'SSG.COM에서 최대 12%', '       스마일클럽으로 5% 적립', '       쓱·새벽·트레이더스 구매 시 (구매 전 스마일클럽 적립 아이콘을 꼭 확인해주세요)', '       SSG.COM 카드 Edition2로 7% 적립', '어디서나 한도없이 0.5%', 'SSG.COM카드 Edition2로 어디서나 한도없이 0.5% 적립 (SSG.COM 외 모든 가맹점)', 'SSG Money 최대 12% 적립', '스마일클럽 5% 적립 + SSG.COM카드 Edition2 최대 7% 적립', '      스마일클럽 5% 적립은 쓱·새벽배송·트레이더스 이용 시에 한함', '      SSG.COM카드 Edition2 최대 7% 적립', '      SSG.COM에서 결제 시 7%(1만 쓱머니 한도), 그 외 가맹점 0.5% 적립(적립한도 제한 없음)', '      무이자 할부 및 현대카드에서 제공하는 다른 할인 서비스 이용 시 적립 제외', '스마일클럽 가입비 3,900원 매월 100% 지원', '      1. 스마일클럽 자동 가입에', '      2. SSGPAY 내 카드 자동 등록', '      3. 스마일클럽 정기결제수단 자동 등록 및 월 이용료 3,900원 지원까지! (단, 해당 카드를 월 정기결제 수단에 등록한 경우에 한함)', 'SSG.COM카드 Edition2를 스마일클럽 월 정기결제 수단에 등록 및 전월 이용금액 30만원 이상 시 혜택 제공', 'TIP. 스마일클럽 가입 시 SSG.COM 혜택', '(SSG 가입 시) 15,000원만 담아도 쓱 · 새벽배송 무료배송', '쓱 · 새벽배송 · 트레이더스 SSG머니 최대 5% 적립', '장바구니 최대 10% 할인쿠폰', '매월 4장씩 최대 12% 할인쿠폰', '스마일배송 1만 5천원 이상 무제한 무료배송', '스마일배송 상품 스마일캐시 1% 적립', '스마일클럽 단독 혜택 스타벅스 상품 전용 딜', '스마일클럽 가입비 매월 3,900원 지원', '월 1회, 매달 스마일클럽 정기결제일에 혜택 제공', 'SSG.COM카드 Edition2는 최초 발급 시, 스마일클럽 월 정기결제 수단에 자동 등록 됨', '전월 이용금액 30만원 미만 시, SSG.COM카드 Edition2로 스마일클럽 정기 결제 금액이 자동 결제됨', '스마일클럽 무료 이용 기간이라면 정기결제 금액 지원 대신 SSG머니 3,900원 제공', '      스타벅스 자동 충전, 생활요금(통신요금, 아파트관리비 등) 정기결제 신청 및 이체 시 최대 1만원 청구 할인', '      2. 스타벅스 자동 충전 또는 생활요금 정기결제 신청(각 항목당 할인한도 5천원, 최대 1만원 할인)', '      정기결제 신청 후 카드 결제일에 따라 매출 발생 다음 달 또는 다다음 달 청구 할인 혜택 적용', '      단, 청구 할인 제공 일정은 당사 또는 신청인 사정에 의해 상이할 수 있음', '      3. 쏘카 1만원 할인쿠폰', '      쿠폰은 등록일 포함 30일간 이용 가능', '실물 SSG.COM카드 Edition2 수령 전 SSGPAY로 결제 시 건당 100만원 이하 결제 가능 *단, 본인 확인(신분증 확인 및 1원 인증) 완료한 경우에 한하며, 건당 100만원 초과 시 실물카드 수령 후 결제 가능', 'SSG머니 최대 적립 12%에서 5%는 SSG.COM에서 제공하는 멤버십 서비스로 SSG.COM 사정에 따라 변경 가능함', '카드 이용대금 연체 시 약정금리 + 연체가산금리 3%의 연체이자율이 적용됩니다. (회원별, 이용 상품별 차등적용 / 법정 최고금리 20% 이내) 단, 연체 발생시점에 약정금리가 없는 경우 아래와 같이 적용', '일시불 : 거래 발생시점 기준 최소 기간 (2개월)의 유이자 할부 약정금리 + 연체가산금리 3%', '무이자할부 : 거래발생시점 기준 동일한 할부 계약 기간의 유이자할부 약정금리 + 연체가산금리 3%']
'@
$ws.Range("G5").Value = $G5

# --- Row 6 ---
$A6 = @'
스마일클럽 웰컴 5천원 쿠폰
'@
$ws.Range("A6").Value = $A6
$B6 = @'
2월의신규가입혜택
'@
$ws.Range("B6").Value = $B6
$C6 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000003992
'@
$ws.Range("C6").Value = $C6
$D6 = @'
[스마일클럽] 2월 매일매일 스마일
'@
$ws.Range("D6").Value = $D6
$ws.Range("E6").NumberFormat = "@"
$E6 = @'
2023-02-01
'@
$ws.Range("E6").Value = $E6
$ws.Range("F6").NumberFormat = "@"
$F6 = @'
2023-02-28
'@
$ws.Range("F6").Value = $F6
$G6 = @'
['이벤트/쿠폰 > [스마일클럽] 2월 매일매일 스마일', '스마일클럽', '[스마일클럽] 2월 매일매일 스마일', ' 지금 스마일클럽 가입하고 모든 혜택 받기 ', '본 쿠폰은 이벤트 기간 내 SSG.COM에서 스마일클럽에 최초 신규 가입한 고객님에 한해 ID당 1회 발급됩니다.', '쿠폰은 가입 차주 금요일 이내 자동 지급되며, 앱푸시나 문자메시지 등을 통해 별도 안내 예정입니다.', '본 쿠폰은 5,100원 이상 구매 시 5,000원 할인됩니다. (할인액 및 배송비 제외한 구매 금액 기준으로 쿠폰 적용됨)', '일부 상품 및 브랜드는 쿠폰 적용 제외될 수 있습니다.', 'SSG.COM의 쿠폰은 결제 시 최적 할인에 의해 자동 설정되며, 직접 변경이 가능합니다.', '본 이벤트는 당사 사정으로 내용이 변경되거나 종료될 수 있습니다.']
'@
$ws.Range("G6").Value = $G6

# --- Row 7 ---
$A7 = @'
첫 구매 고객 스페셜 혜택
'@
$ws.Range("A7").Value = $A7
$B7 = @'
메가박스 영화 관람권 5,900원
'@
$ws.Range("B7").Value = $B7
$C7 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000003994
'@
$ws.Range("C7").Value = $C7
$D7 = @'
[쓱- 첫구매 응원 이벤트] 메가박스 스페셜 할인 혜택
'@
$ws.Range("D7").Value = $D7
$ws.Range("E7").NumberFormat = "@"
$E7 = @'
2023-01-27
'@
$ws.Range("E7").Value = $E7
$ws.Range("F7").NumberFormat = "@"
$F7 = @'
2023-02-09
'@
$ws.Range("F7").Value = $F7
$G7 = @'
['이벤트/쿠폰 > [쓱- 첫구매 응원 이벤트] 메가박스 스페셜 할인 혜택', '스마일클럽', '[쓱- 첫구매 응원 이벤트] 메가박스 스페셜 할인 혜택', '메가박스 스페셜 할인 혜택', '메가박스 스페셜 할인 혜택으로 가족&연인과 즐거운 시간 보내세요!', 'SSG.COM 첫구매 응원 이벤트', '쓱닷컴 첫구매 고객님과 1년만에 다시 오신 고객님을 위한 메가박스 전용 특별 할인쿠폰', '       4,000원 할인쿠폰', '       첫 구매 고객 전용 쿠폰 발급 대상', '       쿠폰발급 및 사용기간', '       2023년 1월 27일(금) 10:00시 부터 선착순 발급 및 발급 후 2월 9일(목)23:59까지 사용가능', '       본 쿠폰은 지정된 메가박스 스페셜 할인 2D 영화 관람권에 한하여 적용 가능하며, 1장 당 상품 1개에 적용됩니다.', '       본 이벤트는 당사 사정에 따라 변동 및 조기 종료될 수 있습니다.', '       정상 판매가 : 13,000원 / SSG.COM 스페셜 할인가 : 9,900원']
'@
$ws.Range("G7").Value = $G7

# --- Row 8 ---
$A8 = @'
폴레드 2/2(목) 11:00AM
'@
$ws.Range("A8").Value = $A8
$B8 = @'
국민육아핫템! 폴레드 최대79% 역대급 할인
'@
$ws.Range("B8").Value = $B8
$C8 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000004147
'@
$ws.Range("C8").Value = $C8
$D8 = @'
폴레드 @SSG.LIVE 2/2(목) 11:00AM
'@
$ws.Range("D8").Value = $D8
$ws.Range("E8").NumberFormat = "@"
$E8 = @'
2023-01-30
'@
$ws.Range("E8").Value = $E8
$ws.Range("F8").NumberFormat = "@"
$F8 = @'
2023-02-02
'@
$ws.Range("F8").Value = $F8
$G8 = @'
['이벤트/쿠폰 > 폴레드 @SSG.LIVE 2/2(목) 11:00AM', '스마일클럽', '                    \xa0SSG.LIVE 사은품 지급 및 이벤트 혜택 당첨 주의사항', '\xa0- 사은품 지급 및 이벤트 혜택 제공', '\xa0- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기']
'@
$ws.Range("G8").Value = $G8

# --- Row 9 ---
$A9 = @'
톰포드뷰티 2/2(목) 7PM
'@
$ws.Range("A9").Value = $A9
$B9 = @'
톰포드뷰티 발렌타인 선물제안 로즈프릭 에디션 긴급공수 & 상품권 증정
'@
$ws.Range("B9").Value = $B9
$C9 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000004083
'@
$ws.Range("C9").Value = $C9
$D9 = @'
톰포드뷰티 @SSG.LIVE 2/2(목) 7PM
'@
$ws.Range("D9").Value = $D9
$ws.Range("E9").NumberFormat = "@"
$E9 = @'
2023-01-26
'@
$ws.Range("E9").Value = $E9
$ws.Range("F9").NumberFormat = "@"
$F9 = @'
2023-02-02
'@
$ws.Range("F9").Value = $F9
$G9 = @'
['이벤트/쿠폰 > 톰포드뷰티 @SSG.LIVE 2/2(목) 7PM', '스마일클럽', '이벤트 혜택 당첨 주의사항', '- 사은품 지급 및 이벤트 혜택 제공', '- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기']
'@
$ws.Range("G9").Value = $G9

# --- Row 10 ---
$A10 = @'
웨스틴조선서울 2/2(목) 8PM
'@
$ws.Range("A10").Value = $A10
$B10 = @'
주중(일-목) 이그제큐티브룸 구매시 클럽라운지 1인 추가 무료
'@
$ws.Range("B10").Value = $B10
$C10 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000003995
'@
$ws.Range("C10").Value = $C10
$D10 = @'
웨스틴조선서울 이그제큐티브 핫딜 29만원대부터~ @SSG.LIVE 2/2(목) 8PM
'@
$ws.Range("D10").Value = $D10
$ws.Range("F10").NumberFormat = "@"
$F10 = @'
2023-02-02
'@
$ws.Range("F10").Value = $F10
$G10 = @'
['이벤트/쿠폰 > 웨스틴조선서울 이그제큐티브 핫딜 29만원대부터~ @SSG.LIVE 2/2(목) 8PM', '스마일클럽', '웨스틴조선서울 이그제큐티브 핫딜 29만원대부터~ @SSG.LIVE 2/2(목) 8PM']
'@
$ws.Range("G10").Value = $G10

# --- Row 11 ---
$A11 = @'
빈폴키즈 2/3(금) 11AM
'@
$ws.Range("A11").Value = $A11
$B11 = @'
우리아이 조카 선물로 딱! / LIVE에서만 책가방 최대 28% 압도적 혜택
'@
$ws.Range("B11").Value = $B11
$C11 = @'
https://event.ssg.com/eventDetail.ssg?nevntId=1000000004084
'@
$ws.Range("C11").Value = $C11
$D11 = @'
빈폴키즈 @SSG.LIVE 2/3(금) 11:00
'@
$ws.Range("D11").Value = $D11
$ws.Range("F11").NumberFormat = "@"
$F11 = @'
2023-02-03
'@
$ws.Range("F11").Value = $F11
$G11 = @'
['이벤트/쿠폰 > 빈폴키즈 @SSG.LIVE 2/3(금) 11:00', '스마일클럽', '- 사은품 지급 및 이벤트 혜택 제공', '- 사은품 지급 및 이벤트 혜택 제공 관련 업무 종료 후 즉시 파기']
'@
$ws.Range("G11").Value = $G11

# Remove rows 12 and 13 (entries no longer present)
$ws.Rows("12:13").Delete()
